$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("C8").Value = "[Lab 5 Sampling Distributions](https://crumplab.github.io/psyc7709Lab/articles/Lab5_Sampling_Distributions.html)"
$ws.Range("C8").RowHeight = 51
$ws.Range("C9").Select() | Out-Null
